$wb = $excel.ActiveWorkbook

# --- Sheet: ИсходныеДанные (source data) - add new header columns AX1:BO1 ---
$wsData = $wb.Worksheets.Item("ИсходныеДанные")
$wsData.Cells.Item(1, 50).Value = "IS_Service_type"
$wsData.Cells.Item(1, 51).Value = "IS_Service_type_Month"
$wsData.Cells.Item(1, 52).Value = "IS_Product_type"
$wsData.Cells.Item(1, 53).Value = "IS_Product_type_Month"
$wsData.Cells.Item(1, 54).Value = "Pdr_Proj"
$wsData.Cells.Item(1, 55).Value = "Pdr_Proj_Month"
$wsData.Cells.Item(1, 56).Value = "Proj_Pdr"
$wsData.Cells.Item(1, 57).Value = "Proj_Pdr_Month"
$wsData.Cells.Item(1, 58).Value = "FN_Month"
$wsData.Cells.Item(1, 59).Value = "UHCost_KV1"
$wsData.Cells.Item(1, 60).Value = "UMCost_KV1"
$wsData.Cells.Item(1, 61).Value = "UHCost_KV2"
$wsData.Cells.Item(1, 62).Value = "UMCost_KV2"
$wsData.Cells.Item(1, 63).Value = "UHCost_KV3"
$wsData.Cells.Item(1, 64).Value = "UMCost_KV3"
$wsData.Cells.Item(1, 65).Value = "UHCost_KV4"
$wsData.Cells.Item(1, 66).Value = "UMCost_KV4"
$wsData.Cells.Item(1, 67).Value = "ISDogName"

# --- Sheet: Настройки (settings) - add new configuration rows 16:34 ---
$wsSettings = $wb.Worksheets.Item("Настройки")
$wsSettings.Cells.Item(16, 1).Value = "IS_Service_type"
$wsSettings.Cells.Item(16, 2).Value = "Тип сервиса (ИСУ, КИС, ЛИС, ПУ, ..)"
$wsSettings.Cells.Item(16, 3).Value = "AX"
$wsSettings.Cells.Item(16, 4).Value = "AY"
$wsSettings.Cells.Item(17, 1).Value = "IS_Product_type"
$wsSettings.Cells.Item(17, 2).Value = "Тип системы (SAP, БК, ЛИМС, MES,…)"
$wsSettings.Cells.Item(17, 3).Value = "AZ"
$wsSettings.Cells.Item(17, 4).Value = "BA"
$wsSettings.Cells.Item(18, 1).Value = "Pdr_Proj"
$wsSettings.Cells.Item(18, 2).Value = "Группировка Подразделение+Проект"
$wsSettings.Cells.Item(18, 3).Value = "BB"
$wsSettings.Cells.Item(18, 4).Value = "BC"
$wsSettings.Cells.Item(19, 1).Value = "Proj_Pdr"
$wsSettings.Cells.Item(19, 2).Value = "Группировка Проект+Подразделение"
$wsSettings.Cells.Item(19, 3).Value = "BD"
$wsSettings.Cells.Item(19, 4).Value = "BE"
$wsSettings.Cells.Item(20, 1).Value = "Portfolio"
$wsSettings.Cells.Item(20, 2).Value = "Портфель проектов"
$wsSettings.Cells.Item(20, 3).Value = "AV"
$wsSettings.Cells.Item(20, 4).Value = "AW"
$wsSettings.Cells.Item(21, 1).Value = "Personal_email"
$wsSettings.Cells.Item(21, 2).Value = "Признак отправлять сообщение лично или в общей массе"
$wsSettings.Cells.Item(21, 3).Value = "AR"
$wsSettings.Cells.Item(22, 1).Value = "user_email"
$wsSettings.Cells.Item(22, 2).Value = "Почтовый адрес пользователя"
$wsSettings.Cells.Item(22, 3).Value = "AS"
$wsSettings.Cells.Item(23, 1).Value = "boss_email"
$wsSettings.Cells.Item(23, 2).Value = "Почтовый адрес руководителя данного пользователя"
$wsSettings.Cells.Item(23, 3).Value = "AT"
$wsSettings.Cells.Item(24, 1).Value = "Contract"
$wsSettings.Cells.Item(24, 2).Value = "Доходный договор"
$wsSettings.Cells.Item(24, 3).Value = "AU"
$wsSettings.Cells.Item(25, 1).Value = "FN"
$wsSettings.Cells.Item(25, 2).Value = "Функциональное направление (или подразделение)"
$wsSettings.Cells.Item(25, 3).Value = "C"
$wsSettings.Cells.Item(25, 4).Value = "BF"
$wsSettings.Cells.Item(26, 1).Value = "UHCost_KV1"
$wsSettings.Cells.Item(26, 2).Value = "Часовая ставка в 1-м квартале"
$wsSettings.Cells.Item(26, 3).Value = "BF"
$wsSettings.Cells.Item(27, 1).Value = "UMCost_KV1"
$wsSettings.Cells.Item(27, 2).Value = "Месячная ставка в 1-м квартале"
$wsSettings.Cells.Item(27, 3).Value = "BG"
$wsSettings.Cells.Item(28, 1).Value = "UHCost_KV2"
$wsSettings.Cells.Item(28, 2).Value = "Часовая ставка во 2-м квартале"
$wsSettings.Cells.Item(28, 3).Value = "BH"
$wsSettings.Cells.Item(29, 1).Value = "UMCost_KV2"
$wsSettings.Cells.Item(29, 2).Value = "Месячная ставка во 2-м квартале"
$wsSettings.Cells.Item(29, 3).Value = "BI"
$wsSettings.Cells.Item(30, 1).Value = "UHCost_KV3"
$wsSettings.Cells.Item(30, 2).Value = "Часовая ставка в 3-м квартале"
$wsSettings.Cells.Item(30, 3).Value = "BJ"
$wsSettings.Cells.Item(31, 1).Value = "UMCost_KV3"
$wsSettings.Cells.Item(31, 2).Value = "Месячная ставка в 3-м квартале"
$wsSettings.Cells.Item(31, 3).Value = "BK"
$wsSettings.Cells.Item(32, 1).Value = "UHCost_KV4"
$wsSettings.Cells.Item(32, 2).Value = "Часовая ставка в 4-м квартале"
$wsSettings.Cells.Item(32, 3).Value = "BL"
$wsSettings.Cells.Item(33, 1).Value = "UMCost_KV4"
$wsSettings.Cells.Item(33, 2).Value = "Месячная ставка в 4-м квартале"
$wsSettings.Cells.Item(33, 3).Value = "BM"
$wsSettings.Cells.Item(34, 1).Value = "ISDogName"
$wsSettings.Cells.Item(34, 2).Value = "Название ИС из договора"
$wsSettings.Cells.Item(34, 3).Value = "BO"
